$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-07 Saturday" "2024-09-08 Sunday"

Replace-Text "18×50=900" "12×91=1092"
Replace-Text "45×95=4275" "47×92=4324"
Replace-Text "78×90=7020" "69×97=6693"
Replace-Text "47×54=2538" "53×20=1060"
Replace-Text "34×38=1292" "31×27=837"

Replace-Text "22×47=1034" "45×21=945"
Replace-Text "98×88=8624" "80×83=6640"
Replace-Text "58×79=4582" "51×59=3009"
Replace-Text "47×65=3055" "98×60=5880"
Replace-Text "11×32=352" "93×96=8928"

Replace-Text "96×27=2592" "64×90=5760"
Replace-Text "89×58=5162" "38×62=2356"
Replace-Text "98×45=4410" "11×32=352"
Replace-Text "11×74=814" "78×91=7098"
Replace-Text "24×27=648" "71×78=5538"

Replace-Text "37×74=2738" "11×13=143"
Replace-Text "67×43=2881" "34×46=1564"
Replace-Text "57×40=2280" "57×12=684"
Replace-Text "96×12=1152" "67×33=2211"
Replace-Text "54×79=4266" "53×79=4187"

Replace-Text "26×77=2002" "64×13=832"
Replace-Text "34×63=2142" "94×46=4324"
Replace-Text "55×76=4180" "96×97=9312"
Replace-Text "93×83=7719" "23×91=2093"
Replace-Text "29×26=754" "24×95=2280"
